# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 25 de Mayo de 2020 a las 22:05"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1700877
$ws.Range("C4").Value = 14441
$ws.Range("D4").Value = 457369
$ws.Range("E4").Value = 1143823
$ws.Range("G4").Value = 385
$ws.Range("H4").Value = 99685

# Row 11: Alemania
$ws.Range("B11").Value = 180789
$ws.Range("C11").Value = 461
$ws.Range("E11").Value = 11161
$ws.Range("G11").Value = 57
$ws.Range("H11").Value = 8428

# Row 16: Canada
$ws.Range("D16").Value = 44530
$ws.Range("E16").Value = 34609

# Row 25: Ecuador
$ws.Range("D25").Value = 18003
$ws.Range("E25").Value = 16149

# Row 190: Gambia
$ws.Range("D190").Value = 17
$ws.Range("E190").Value = 7
